# Update crypto price/volume figures (GitHub Actions scheduled refresh).
# Values in columns D (Price) and E (Volume(1h)) are stored as TEXT in the
# workbook (even though many look numeric), so every .Value assignment below
# is apostrophe-prefixed to force Excel to keep them as text instead of
# auto-converting look-alike numbers (e.g. "69.37", "8.21") into the Number
# type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.Value = "'" + $Text
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.655.96"
Set-TextValue $ws.Range("E2") "  +1.36%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.985.78"

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "379.87"
Set-TextValue $ws.Range("E5") "  +3.06%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "104.58"
Set-TextValue $ws.Range("E6") "  +1.53%  "

# Row 7 - XRP
Set-TextValue $ws.Range("E7") "  +0.87%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  -0.01%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.596"
Set-TextValue $ws.Range("E9") "  +2.09%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "37.30"
Set-TextValue $ws.Range("E10") "  +2.62%  "

# Row 11 - TRON
Set-TextValue $ws.Range("E11") "  +0.18%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("E12") "  +1.97%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "3.453.98"
Set-TextValue $ws.Range("E13") "  +2.96%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "18.45"
Set-TextValue $ws.Range("E14") "  +1.41%  "

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") "7.61"
Set-TextValue $ws.Range("E15") "  +3.62%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.990.81"
Set-TextValue $ws.Range("E16") "  +2.98%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.980"
Set-TextValue $ws.Range("E17") "  +6.13%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "51.596.96"

# Row 19 - ImmutableX
Set-TextValue $ws.Range("D19") "3.34"
Set-TextValue $ws.Range("E19") "  +3.93%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("E20") "  +4.12%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.99"
Set-TextValue $ws.Range("E21") "  +1.26%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("D22") "0.0₃0965"
Set-TextValue $ws.Range("E22") "  +2.82%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "69.37"
Set-TextValue $ws.Range("E23") "  +1.78%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "262.41"
Set-TextValue $ws.Range("E24") "  +1.61%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("E25") "  +7.05%  "

# Row 26 - Filecoin
Set-TextValue $ws.Range("D26") "8.21"
Set-TextValue $ws.Range("E26") "  +18.96%  "

# Row 27 - RenderToken
Set-TextValue $ws.Range("D27") "7.73"
Set-TextValue $ws.Range("E27") "  +24.64%  "

# Row 28 - Hedera
Set-TextValue $ws.Range("D28") "0.116"
Set-TextValue $ws.Range("E28") "  +14.34%  "

# Row 29 - Kaspa
Set-TextValue $ws.Range("E29") "  -1.35%  "

# Row 30 - Dai
Set-TextValue $ws.Range("E30") "  +0.06%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Range("D31") "25.94"
Set-TextValue $ws.Range("E31") "  +1.63%  "

# Row 32 - Cosmos
Set-TextValue $ws.Range("D32") "9.90"
Set-TextValue $ws.Range("E32") "  +0.49%  "

# Row 33 - InjectiveProtocol
Set-TextValue $ws.Range("D33") "35.16"
Set-TextValue $ws.Range("E33") "  +2.31%  "

# Row 34 / Row 35 - Toncoin and OKB swap positions (OKB now ranks above Toncoin)
Set-TextValue $ws.Range("B34") "OKB"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "51.16"
Set-TextValue $ws.Range("E34") "  +0.73%  "

Set-TextValue $ws.Range("B35") "Toncoin"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D35") "2.08"
Set-TextValue $ws.Range("E35") "  -1.87%  "

# Row 36 - VeChain
Set-TextValue $ws.Range("D36") "0.0448"
Set-TextValue $ws.Range("E36") "  +6.57%  "

# Row 37 - FirstDigitalUSD
Set-TextValue $ws.Range("E37") "  +0.09%  "

# Row 38 - LidoDAOToken
Set-TextValue $ws.Range("E38") "  +2.03%  "

# Row 39 - Celestia
Set-TextValue $ws.Range("D39") "17.16"
Set-TextValue $ws.Range("E39") "  +1.25%  "

# Row 40 - Stacks
Set-TextValue $ws.Range("D40") "2.60"
Set-TextValue $ws.Range("E40") "  -0.78%  "

# Row 41 - ARBITRUM
Set-TextValue $ws.Range("E41") "  +1.28%  "

# Row 42 - Stellar
Set-TextValue $ws.Range("E42") "  +3.04%  "

# Row 43 - Monero
Set-TextValue $ws.Range("D43") "125.61"
Set-TextValue $ws.Range("E43") "  +6.31%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "21.77"
Set-TextValue $ws.Range("E44") "  -1.38%  "

# Row 45 - TheGraph
Set-TextValue $ws.Range("D45") "0.284"
Set-TextValue $ws.Range("E45") "  +19.89%  "

# Row 46 - WEMIXToken
Set-TextValue $ws.Range("E46") "  -1.13%  "

# Row 47 - ApeXProtocol
Set-TextValue $ws.Range("D47") "2.38"
Set-TextValue $ws.Range("E47") "  +2.90%  "

# Row 48 - Maker
Set-TextValue $ws.Range("D48") "2.038.70"
Set-TextValue $ws.Range("E48") "  +1.52%  "

# Row 49 - NEARProtocol
Set-TextValue $ws.Range("E49") "  +3.49%  "

# Row 50 - BEAM
Set-TextValue $ws.Range("E50") "  +8.98%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "58.43"
Set-TextValue $ws.Range("E51") "  +3.62%  "
